# Update version string across the workbook for release
# "mines - version 1.0.0 (Feb 3 2026)" built on "February 03 2026 10.14.00 EST"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# Update the "Version:" line on the About sheet
$wsAbout.Range("A2").Value = "Version: $newVersion"

# Update the "Recommended Citation:" line on the About sheet
$oldCitation = $wsAbout.Range("A6").Value2
$newCitation = $oldCitation -replace [regex]::Escape($oldVersion), $newVersion
$wsAbout.Range("A6").Value = $newCitation

# Update the build_version column (S2:S14) on the data sheet
for ($r = 2; $r -le 14; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
